$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.317.80"
$ws.Range("E2").Value = "  +3.18%  "
$ws.Range("D3").Value = "1.987.49"
$ws.Range("E3").Value = "  +6.32%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'0.7905"
$ws.Range("E5").Value = "  +68.23%  "
$ws.Range("D6").Value = "'253.00"
$ws.Range("E6").Value = "  +3.88%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'0.3395"
$ws.Range("E8").Value = "  +18.13%  "
$ws.Range("D9").Value = "'25.78"
$ws.Range("E9").Value = "  +17.08%  "
$ws.Range("D10").Value = "'0.06947"
$ws.Range("E10").Value = "  +7.61%  "
$ws.Range("D11").Value = "'0.8372"
$ws.Range("E11").Value = "  +16.16%  "
$ws.Range("D12").Value = "'0.08112"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.992.52"
$ws.Range("E13").Value = "  +6.44%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'100.63"
$ws.Range("E14").Value = "  +5.06%  "
$ws.Range("D15").Value = "'5.444"
$ws.Range("E15").Value = "  +6.22%  "
$ws.Range("D16").Value = "'271.98"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").Value = "31.347.88"
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("D18").Value = "'13.91"
$ws.Range("E18").Value = "  +7.23%  "
$ws.Range("D19").Value = "'0.000007933"
$ws.Range("E19").Value = "  +5.43%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.252.40"
$ws.Range("E20").Value = "  +6.43%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.715"
$ws.Range("E21").Value = "  +9.29%  "
$ws.Range("D22").Value = "'1.006"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").Value = "'6.944"
$ws.Range("E24").Value = "  +11.45%  "
$ws.Range("D25").Value = "'9.647"
$ws.Range("E25").Value = "  +6.78%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'165.05"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1476"
$ws.Range("E27").Value = "  +53.84%  "
$ws.Range("D28").Value = "'19.71"
$ws.Range("E28").Value = "  +5.63%  "
$ws.Range("D29").Value = "'2.172"
$ws.Range("E29").Value = "  +15.92%  "
$ws.Range("E30").Value = "  +6.72%  "
$ws.Range("D31").Value = "'1.360"
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("D32").Value = "'4.561"
$ws.Range("E32").Value = "  +8.44%  "
$ws.Range("D33").Value = "'4.319"
$ws.Range("E33").Value = "  +5.43%  "
$ws.Range("D34").Value = "'0.05177"
$ws.Range("E34").Value = "  +7.60%  "
$ws.Range("E35").Value = "  +8.49%  "
$ws.Range("D36").Value = "'0.7512"
$ws.Range("E36").Value = "  +9.11%  "
$ws.Range("D37").Value = "'2.799"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "'1.004"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01999"
$ws.Range("E39").Value = "  +6.54%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.924"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.627"
$ws.Range("E41").Value = "  +6.70%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'78.40"
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4639"
$ws.Range("E43").Value = "  +10.02%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.063"
$ws.Range("E44").Value = "  +6.48%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'105.46"
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'0.8519"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.02"
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.479"
$ws.Range("E49").Value = "  +7.71%  "
$ws.Range("D50").Value = "'0.4277"
$ws.Range("E50").Value = "  +9.17%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'36.48"
$ws.Range("E51").Value = "  +3.58%  "
